$d = $word.ActiveDocument

# Mark every inline picture's run as "no proofing" (w:rPr/w:noProof),
# matching what Word stamps on runs containing drawings it rendered.
foreach ($shape in $d.InlineShapes) {
    $shape.Range.NoProofing = $true
}

# The last (previously empty) paragraph before the sectPr now contains
# a single space character.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter(" ")

Write-Output "Applied noProof to $($d.InlineShapes.Count) inline shapes; appended space to last paragraph."
